$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control_obs")
$ws.Range("AH20").Formula = "=SUM(AH2:AH18)"
$ws.Range("AG20").Copy()
$ws.Range("AH20").PasteSpecial(-4122)
Write-Output "done"
